$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3 (SamplesTab): replace the samples query with the reduced-column version
$newSamplesQuery = @'
SELECT
    DISTINCT (smp.sample_id) AS "Sample ID",
    sp.participant_id AS "Participant ID", 
    s.study_name AS "Study Name",
    s.phs_accession AS Accession
FROM 
    df_participant sp
JOIN 
    df_study s ON sp."study.phs_accession" = s.phs_accession
JOIN 
    df_sample smp ON smp."participant.study_participant_id" = sp.study_participant_id
JOIN
    df_diagnosis d ON d."participant.study_participant_id" = sp.study_participant_id
JOIN
    df_program p ON p.program_acronym = s."program.program_acronym"
JOIN
    df_file f1 ON f1."sample.sample_id" = smp.sample_id
JOIN
    df_genomic_info gi ON gi."file.file_id" = f1.file_id
WHERE 
    s.phs_accession = 'phs001437' AND gi.platform = 'Illumina'
ORDER BY 
    smp.sample_id ASC
LIMIT 100;
'@
$ws.Range("B3").Value = $newSamplesQuery

# Rows 3 and 4 no longer carry the TSV/Web output-file columns (D, E)
$ws.Range("D3:E4").ClearContents()

# Update the active selection / scroll position left by the author
$ws.Range("B3").Select()
